$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneratorsOfNode")

# Remove the obsolete/duplicate "Wind offshore grounded"/"Wind offshore
# floating" rows that used stale ASCII node names (Helgolander Bucht,
# Nordsoen, Sorlige Nordsjo I/II) instead of the canonical OffshoreNodes
# spellings. Two separate 4-row blocks are removed (rows 641-644 first, so
# the still-unshifted row numbers for the first block, 623-626, stay valid),
# shifting everything below each deletion up; dimension goes from
# A1:B741 to A1:B733 (8 rows removed total).
$ws.Range("A641:B644").EntireRow.Delete()
$ws.Range("A623:B626").EntireRow.Delete()
